# Edit: Add an "output_name" column (placed immediately before the existing
# "id" column) to every worksheet that represents a LinkML "Source" class,
# i.e. every sheet whose header row ends in input/output/type/id.
#
# Commit message: "Remove name from sequence, and add output_name to source (#11)"
# (the "name" removal does not affect this workbook's generated sheets, so
# only the "output_name" addition is applied here.)

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> column letter that currently holds "id" (i.e. the
# last column of the header row). "output_name" is written into that column
# and "id" is moved one column to the right.
$sheetsToUpdate = @{
    "ManuallyTypedSource"              = "H"
    "UploadedFileSource"               = "H"
    "RepositoryIdSource"               = "F"
    "AddGeneIdSource"                  = "H"
    "GenomeCoordinatesSource"          = "K"
    "SequenceCutSource"                = "F"
    "RestrictionEnzymeDigestionSource" = "F"
    "AssemblySource"                   = "F"
    "PCRSource"                        = "H"
    "LigationSource"                   = "F"
    "HomologousRecombinationSource"    = "F"
    "GibsonAssemblySource"             = "F"
    "RestrictionAndLigationSource"     = "G"
    "CRISPRSource"                     = "G"
    "OligoHybridizationSource"         = "G"
    "PolymeraseExtensionSource"        = "D"
}

foreach ($sheetName in $sheetsToUpdate.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $idColLetter = $sheetsToUpdate[$sheetName]
    $idCol = $ws.Range($idColLetter + "1").Column
    $newIdCol = $idCol + 1

    # Shift "id" one column to the right, then write "output_name" into the
    # column that used to hold "id".
    $ws.Cells.Item(1, $newIdCol).Value = "id"
    $ws.Cells.Item(1, $idCol).Value = "output_name"
}
